$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Split the single long "customer feedback" paragraph (old row 9) into
#    four separate rows, each holding one sentence, by inserting three new
#    rows right after the existing one.
# ---------------------------------------------------------------------------
$ws.Rows("10:12").Insert()

$ws.Range("A9").Value = "We want to confirm that your shipment of supplied above products have met our expectations completely. "
$ws.Range("A10").Value = "Product quality and performance of the product found best along with on time response and delivery time from "
$ws.Range("A11").Value = "your company."
$ws.Range("A12").Value = "Furthermore, we expect to continue good results with application of your manufactured products on our oil fields."

Write-Host "text set"
